# Insert two new data rows at row 820 (pushing the existing rows 820-920 down
# to 822-922, with all their data unchanged), then populate the two new rows
# with the new weekly price observations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 820 twice (Excel copies the
# formatting, including the date number format, from the row above on insert).
$ws.Rows.Item(820).Insert()
$ws.Rows.Item(820).Insert()

# New row 820: Primera, $/caja 60 unidades, Region de Arica y Parinacota
$ws.Cells.Item(820, 1).Value = 8
$ws.Cells.Item(820, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(820, 3).Value = "Coquimbo"
$ws.Cells.Item(820, 4).Value = 45124
$ws.Cells.Item(820, 5).Value = 4
$ws.Cells.Item(820, 6).Value = 100112043
$ws.Cells.Item(820, 7).Value = "Pepino ensalada"
$ws.Cells.Item(820, 8).Value = "Sin especificar"
$ws.Cells.Item(820, 9).Value = "Primera"
$ws.Cells.Item(820, 10).Value = 560
$ws.Cells.Item(820, 11).Value = 10500
$ws.Cells.Item(820, 12).Value = 11000
$ws.Cells.Item(820, 13).Value = 10750
$ws.Cells.Item(820, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(820, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(820, 16).Value = 179
$ws.Cells.Item(820, 17).Value = 60
$ws.Cells.Item(820, 18).Value = "Hortaliza"

# New row 821: Segunda, $/caja 80 unidades, Region de Arica y Parinacota
$ws.Cells.Item(821, 1).Value = 8
$ws.Cells.Item(821, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(821, 3).Value = "Coquimbo"
$ws.Cells.Item(821, 4).Value = 45124
$ws.Cells.Item(821, 5).Value = 4
$ws.Cells.Item(821, 6).Value = 100112043
$ws.Cells.Item(821, 7).Value = "Pepino ensalada"
$ws.Cells.Item(821, 8).Value = "Sin especificar"
$ws.Cells.Item(821, 9).Value = "Segunda"
$ws.Cells.Item(821, 10).Value = 360
$ws.Cells.Item(821, 11).Value = 7500
$ws.Cells.Item(821, 12).Value = 8000
$ws.Cells.Item(821, 13).Value = 7750
$ws.Cells.Item(821, 14).Value = "`$/caja 80 unidades"
$ws.Cells.Item(821, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(821, 16).Value = 97
$ws.Cells.Item(821, 17).Value = 80
$ws.Cells.Item(821, 18).Value = "Hortaliza"
